$d = $word.ActiveDocument

# 1. Delete the block of Q&A paragraphs (and the two blank paragraphs that
#    followed them) right after the opening "NEID looks at the sun..." text.
#    These are paragraphs 3 through 7 (1-indexed):
#      "Q: is the actually the number of photons?"
#      "Q: Why is the range of wavelengths so small?"
#      "Q: Is this averaging over the sun's surface?"
#      (empty)
#      (empty)
$start = $d.Paragraphs(3).Range.Start
$end = $d.Paragraphs(7).Range.End
$d.Range($start, $end).Delete()

# 2. Fill in the empty bullet under "Description of the inputs..." with the
#    new sentence about the L1 raw spectral data.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Description of the inputs*") {
        $newPara = $d.Paragraphs($i + 1)
        $newPara.Range.InsertBefore("L1 raw spectral data of the Sun, which has ~40k fits files.")
        # Restore the Roboto/Times-New-Roman run formatting that matches the
        # rest of the bulleted list (plain InsertBefore resets to defaults).
        $newPara.Range.Font.Name = "Roboto"
        $newPara.Range.Font.NameFarEast = "Times New Roman"
        $newPara.Range.Font.NameBi = "Times New Roman"
        $newPara.Range.Font.Color = 2236962
        break
    }
}

# 3. Rework the "We will compare..." testing-plan sentence.
$d.Content.Find.Execute("vs the number of lines", $true, $false, $false, $false, $false, $true, 1, $false, "versus the number of lines in the catalog that", 2) | Out-Null
$d.Content.Find.Execute(" The one free parameter for this will be N, the number of lines to find. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 4. Replace "Seems to be noisy" with the new wording, then delete the
#    "July 31st" paragraph that directly followed it.
$d.Content.Find.Execute("Seems to be noisy", $true, $false, $false, $false, $false, $true, 1, $false, "There are some fits files where it just seems like background noise, almost as if the detector had a covering over it.", 2) | Out-Null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "July 31*") {
        $d.Paragraphs($i).Range.Delete()
        break
    }
}

# 5. Replace "False dips, because the measure of a dip is inaccurate" with
#    the new wording.
$d.Content.Find.Execute("False dips, because the measure of a dip is inaccurate", $true, $false, $false, $false, $false, $true, 1, $false, "The algorithm may find “false dips”", 2) | Out-Null

# 6. Remove the trailing empty bullet after "...numpy to do math".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*to do math*") {
        $d.Paragraphs($i + 1).Range.Delete()
        break
    }
}
